$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1111197.8
$ws.Range("I2").Value = 1587368.2
$ws.Range("K2").Value = 1587368.2
$ws.Range("M2").Value = -1587255.2
$ws.Range("H3").Value = 48199.6
$ws.Range("J3").Value = 48199.6
$ws.Range("L3").Value = 48199.6
$ws.Range("N3").Value = -48427.6
$ws.Range("H40").Value = 8052.6665
$ws.Range("I40").Value = 9661.666999999999
$ws.Range("J40").Value = 1616.6666
$ws.Range("K40").Value = 9661.666999999999
$ws.Range("L40").Value = 1616.6666
$ws.Range("M40").Value = -9486.666999999999
$ws.Range("N40").Value = -1966.6666
$ws.Range("H48").Value = 4224
$ws.Range("I48").Value = 3439.4
$ws.Range("J48").Value = 5531.6665
$ws.Range("K48").Value = 10318.2
$ws.Range("L48").Value = 16594.9995
$ws.Range("M48").Value = -10026.2
$ws.Range("N48").Value = -17178.9995
$ws.Range("H56").Value = 4224
$ws.Range("I56").Value = 3439.4
$ws.Range("J56").Value = 5531.6665
$ws.Range("K56").Value = 10318.2
$ws.Range("L56").Value = 16594.9995
$ws.Range("M56").Value = -9784.200000000001
$ws.Range("N56").Value = -17662.9995
$ws.Range("H62").Value = 2153.111
$ws.Range("I62").Value = 1625.4286
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 1625.4286
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -1001.4286
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 2153.111
$ws.Range("I65").Value = 1625.4286
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 8127.143
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -5007.143
$ws.Range("N65").Value = -26240
$ws.Range("H69").Value = 4525.423
$ws.Range("I69").Value = 3913
$ws.Range("K69").Value = 11739
$ws.Range("M69").Value = -10865
$ws.Range("H72").Value = 4525.423
$ws.Range("I72").Value = 3913
$ws.Range("K72").Value = 35217
$ws.Range("M72").Value = -30849
$ws.Range("H92").Value = 346.3143
$ws.Range("I92").Value = 234.6923
$ws.Range("K92").Value = 234.6923
$ws.Range("M92").Value = 1013.3077
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H96").Value = 1719
$ws.Range("I96").Value = 2356.6667
$ws.Range("J96").Value = 953.8
$ws.Range("K96").Value = 7070.000100000001
$ws.Range("L96").Value = 2861.4
$ws.Range("M96").Value = -5697.000100000001
$ws.Range("N96").Value = -5607.4
$ws.Range("H102").Value = 48199.6
$ws.Range("J102").Value = 48199.6
$ws.Range("L102").Value = 48199.6
$ws.Range("N102").Value = -54689.6
$ws.Range("H131").Value = 1296.875
$ws.Range("I131").Value = 562.5
$ws.Range("J131").Value = 3500
$ws.Range("K131").Value = 1687.5
$ws.Range("L131").Value = 10500
$ws.Range("M131").Value = 3352.5
$ws.Range("N131").Value = -20580
$ws.Range("H135").Value = 3675.4375
$ws.Range("I135").Value = 5347.4
$ws.Range("J135").Value = 888.8333
$ws.Range("K135").Value = 48126.6
$ws.Range("L135").Value = 7999.4997
$ws.Range("M135").Value = -45591.6
$ws.Range("N135").Value = -13069.4997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10240.2
$ws.Range("I28").Value = 5066.8887
$ws.Range("J28").Value = 56800
$ws.Range("K28").Value = 5066.8887
$ws.Range("L28").Value = 56800
$ws.Range("M28").Value = -4874.8887
$ws.Range("N28").Value = -57184
$ws.Range("H32").Value = 5620.2764
$ws.Range("I32").Value = 2462.111
$ws.Range("K32").Value = 2462.111
$ws.Range("M32").Value = -2175.111
$ws.Range("H61").Value = 1485.5358
$ws.Range("I61").Value = 1369.0385
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1369.0385
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1157.0385
$ws.Range("N61").Value = -3424
$ws.Range("H99").Value = 10240.2
$ws.Range("I99").Value = 5066.8887
$ws.Range("J99").Value = 56800
$ws.Range("K99").Value = 5066.8887
$ws.Range("L99").Value = 56800
$ws.Range("M99").Value = -2071.8887
$ws.Range("N99").Value = -62790
$ws.Range("H101").Value = 71250
$ws.Range("J101").Value = 71250
$ws.Range("L101").Value = 71250
$ws.Range("N101").Value = -77740
$ws.Range("H102").Value = 1233.3684
$ws.Range("I102").Value = 1233.3684
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1233.3684
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 388.6315999999999
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 44900
$ws.Range("J105").Value = 44900
$ws.Range("L105").Value = 44900
$ws.Range("N105").Value = -51888
$ws.Range("H136").Value = 1485.5358
$ws.Range("I136").Value = 1369.0385
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4107.1155
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1557.1155
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 24877.25
$ws.Range("J19").Value = 19836.666
$ws.Range("L19").Value = 19836.666
$ws.Range("N19").Value = -20182.666
$ws.Range("H75").Value = 33585.6
$ws.Range("I75").Value = 7142.6665
$ws.Range("J75").Value = 73250
$ws.Range("K75").Value = 7142.6665
$ws.Range("L75").Value = 73250
$ws.Range("M75").Value = -6206.6665
$ws.Range("N75").Value = -75122
$ws.Range("H78").Value = 33585.6
$ws.Range("I78").Value = 7142.6665
$ws.Range("J78").Value = 73250
$ws.Range("K78").Value = 21427.9995
$ws.Range("L78").Value = 219750
$ws.Range("M78").Value = -16747.9995
$ws.Range("N78").Value = -229110
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H96").Value = 26912.5
$ws.Range("I96").Value = 7533.3335
$ws.Range("J96").Value = 38540
$ws.Range("K96").Value = 7533.3335
$ws.Range("L96").Value = 38540
$ws.Range("M96").Value = -4787.3335
$ws.Range("N96").Value = -44032
$ws.Range("H101").Value = 74400
$ws.Range("J101").Value = 74400
$ws.Range("L101").Value = 74400
$ws.Range("N101").Value = -80890
$ws.Range("H103").Value = 32657
$ws.Range("J103").Value = 32657
$ws.Range("L103").Value = 32657
$ws.Range("N103").Value = -35001
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 989454.2
$ws.Range("I105").Value = 1516349.8
$ws.Range("J105").Value = 1525
$ws.Range("K105").Value = 1516349.8
$ws.Range("L105").Value = 1525
$ws.Range("M105").Value = -1514602.8
$ws.Range("N105").Value = -5019

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 48210.21
$ws.Range("J43").Value = 48210.21
$ws.Range("L43").Value = 48210.21
$ws.Range("N43").Value = -48578.21
$ws.Range("H87").Value = 29996.666
$ws.Range("J87").Value = 29996.666
$ws.Range("L87").Value = 29996.666
$ws.Range("N87").Value = -32368.666
$ws.Range("H90").Value = 29996.666
$ws.Range("J90").Value = 29996.666
$ws.Range("L90").Value = 89989.99800000001
$ws.Range("N90").Value = -101845.998
$ws.Range("H101").Value = 48210.21
$ws.Range("J101").Value = 48210.21
$ws.Range("L101").Value = 48210.21
$ws.Range("N101").Value = -54700.21
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H127").Value = 40347.25
$ws.Range("J127").Value = 40347.25
$ws.Range("L127").Value = 40347.25
$ws.Range("N127").Value = -50267.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 20750
$ws.Range("I99").Value = 16666.666
$ws.Range("J99").Value = 33000
$ws.Range("K99").Value = 16666.666
$ws.Range("L99").Value = 33000
$ws.Range("M99").Value = -13671.666
$ws.Range("N99").Value = -38990
$ws.Range("H100").Value = 1988.8889
$ws.Range("I100").Value = 1700
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1700
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1159
$ws.Range("N100").Value = -4082
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 75800
$ws.Range("J102").Value = 75800
$ws.Range("L102").Value = 75800
$ws.Range("N102").Value = -82290
$ws.Range("H103").Value = 48000
$ws.Range("J103").Value = 48000
$ws.Range("L103").Value = 48000
$ws.Range("N103").Value = -50344
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 4057.1428
$ws.Range("I122").Value = 4045.4546
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 12136.3638
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -9686.363799999999
$ws.Range("N122").Value = -17200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H104").Value = 42870
$ws.Range("J104").Value = 42870
$ws.Range("L104").Value = 42870
$ws.Range("N104").Value = -49858
$ws.Range("H105").Value = 50615
$ws.Range("J105").Value = 50615
$ws.Range("L105").Value = 50615
$ws.Range("N105").Value = -57603
